$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 19 de Octubre de 2020 a las 18:15'
$ws.Cells.Item(4, 2).Value = 8393773
$ws.Cells.Item(4, 3).Value = 5974
$ws.Cells.Item(4, 4).Value = 5463410
$ws.Cells.Item(4, 5).Value = 2705539
$ws.Cells.Item(4, 7).Value = 94
$ws.Cells.Item(4, 8).Value = 224824
$ws.Cells.Item(5, 2).Value = 7574167
$ws.Cells.Item(5, 3).Value = 25929
$ws.Cells.Item(5, 4).Value = 6693491
$ws.Cells.Item(5, 5).Value = 765774
$ws.Cells.Item(5, 7).Value = 260
$ws.Cells.Item(5, 8).Value = 114902
$ws.Cells.Item(14, 2).Value = 741212
$ws.Cells.Item(14, 3).Value = 18804
$ws.Cells.Item(14, 7).Value = 80
$ws.Cells.Item(14, 8).Value = 43726
$ws.Cells.Item(17, 2).Value = 493305
$ws.Cells.Item(17, 3).Value = 1545
$ws.Cells.Item(17, 4).Value = 465021
$ws.Cells.Item(17, 5).Value = 14608
$ws.Cells.Item(17, 7).Value = 41
$ws.Cells.Item(17, 8).Value = 13676
$ws.Cells.Item(19, 2).Value = 423578
$ws.Cells.Item(19, 3).Value = 9338
$ws.Cells.Item(19, 4).Value = 252959
$ws.Cells.Item(19, 5).Value = 134003
$ws.Cells.Item(19, 7).Value = 73
$ws.Cells.Item(19, 8).Value = 36616
$ws.Cells.Item(21, 2).Value = 371154
$ws.Cells.Item(21, 3).Value = 4173
$ws.Cells.Item(21, 5).Value = 69372
$ws.Cells.Item(21, 7).Value = 16
$ws.Cells.Item(21, 8).Value = 9882
$ws.Cells.Item(27, 2).Value = 304367
$ws.Cells.Item(27, 3).Value = 1258
$ws.Cells.Item(27, 4).Value = 272252
$ws.Cells.Item(27, 5).Value = 29867
$ws.Cells.Item(27, 7).Value = 39
$ws.Cells.Item(27, 8).Value = 2248
$ws.Cells.Item(34, 2).Value = 177974
$ws.Cells.Item(34, 3).Value = 4089
$ws.Cells.Item(34, 4).Value = 72553
$ws.Cells.Item(34, 5).Value = 103920
$ws.Cells.Item(34, 7).Value = 79
$ws.Cells.Item(34, 8).Value = 1501
$ws.Cells.Item(41, 2).Value = 121667
$ws.Cells.Item(41, 3).Value = 320
$ws.Cells.Item(41, 4).Value = 98880
$ws.Cells.Item(41, 5).Value = 20584
$ws.Cells.Item(41, 7).Value = 4
$ws.Cells.Item(41, 8).Value = 2203
$ws.Cells.Item(65, 4).Value = 57819
$ws.Cells.Item(65, 5).Value = 68
$ws.Cells.Item(78, 2).Value = 38937
$ws.Cells.Item(78, 3).Value = 1364
$ws.Cells.Item(78, 4).Value = 7006
$ws.Cells.Item(78, 5).Value = 31551
$ws.Cells.Item(78, 7).Value = 35
$ws.Cells.Item(78, 8).Value = 380
$ws.Cells.Item(88, 2).Value = 25802
$ws.Cells.Item(88, 3).Value = 432
$ws.Cells.Item(88, 5).Value = 15293
$ws.Cells.Item(88, 7).Value = 11
$ws.Cells.Item(88, 8).Value = 520
$ws.Cells.Item(90, 2).Value = 23788
$ws.Cells.Item(90, 3).Value = 160
$ws.Cells.Item(90, 4).Value = 17392
$ws.Cells.Item(90, 5).Value = 5550
$ws.Cells.Item(90, 7).Value = 12
$ws.Cells.Item(90, 8).Value = 846
$ws.Cells.Item(95, 2).Value = 17350
$ws.Cells.Item(95, 3).Value = 295
$ws.Cells.Item(95, 4).Value = 10167
$ws.Cells.Item(95, 5).Value = 6729
$ws.Cells.Item(95, 7).Value = 3
$ws.Cells.Item(95, 8).Value = 454
$ws.Cells.Item(99, 2).Value = 15760
$ws.Cells.Item(99, 3).Value = 145
$ws.Cells.Item(99, 4).Value = 11288
$ws.Cells.Item(99, 5).Value = 4232
$ws.Cells.Item(99, 7).Value = 4
$ws.Cells.Item(99, 8).Value = 240
$ws.Cells.Item(108, 2).Value = 11010
$ws.Cells.Item(108, 3).Value = 122
$ws.Cells.Item(108, 4).Value = 8471
$ws.Cells.Item(108, 5).Value = 2404
$ws.Cells.Item(108, 7).Value = 2
$ws.Cells.Item(108, 8).Value = 135
$ws.Cells.Item(112, 2).Value = 10268
$ws.Cells.Item(112, 3).Value = 25
$ws.Cells.Item(112, 4).Value = 9995
$ws.Cells.Item(112, 5).Value = 204
$ws.Cells.Item(123, 2).Value = 5860
$ws.Cells.Item(123, 3).Value = 3
$ws.Cells.Item(123, 4).Value = 4757
$ws.Cells.Item(123, 5).Value = 922
$ws.Cells.Item(125, 2).Value = 5773
$ws.Cells.Item(125, 3).Value = 70
$ws.Cells.Item(125, 4).Value = 3339
$ws.Cells.Item(125, 5).Value = 2311
$ws.Cells.Item(125, 7).Value = 1
$ws.Cells.Item(125, 8).Value = 123
$ws.Cells.Item(137, 1).Value = 'Reunion'
$ws.Cells.Item(137, 2).Value = 4921
$ws.Cells.Item(137, 3).Value = 145
$ws.Cells.Item(137, 4).Value = 4445
$ws.Cells.Item(137, 5).Value = 459
$ws.Cells.Item(137, 8).Value = 17
$ws.Cells.Item(138, 1).Value = 'Republica de Africa Central'
$ws.Cells.Item(138, 2).Value = 4855
$ws.Cells.Item(138, 4).Value = 1924
$ws.Cells.Item(138, 5).Value = 2869
$ws.Cells.Item(138, 8).Value = 62
$ws.Cells.Item(149, 1).Value = 'Principado de Andorra'
$ws.Cells.Item(149, 2).Value = 3623
$ws.Cells.Item(149, 3).Value = 246
$ws.Cells.Item(149, 4).Value = 2273
$ws.Cells.Item(149, 5).Value = 1288
$ws.Cells.Item(149, 7).Value = 3
$ws.Cells.Item(149, 8).Value = 62
$ws.Cells.Item(150, 1).Value = 'Letonia'
$ws.Cells.Item(150, 2).Value = 3494
$ws.Cells.Item(150, 3).Value = 44
$ws.Cells.Item(150, 4).Value = 1341
$ws.Cells.Item(150, 5).Value = 2109
$ws.Cells.Item(150, 8).Value = 44
$ws.Cells.Item(151, 1).Value = 'Mali'
$ws.Cells.Item(151, 2).Value = 3388
$ws.Cells.Item(151, 4).Value = 2586
$ws.Cells.Item(151, 5).Value = 670
$ws.Cells.Item(151, 8).Value = 132
$ws.Cells.Item(152, 2).Value = 2847
$ws.Cells.Item(152, 3).Value = 5
$ws.Cells.Item(152, 5).Value = 1502

Write-Host "Applied updates to paises.xlsx"
